$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.458.09"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").Value = "2.331.17"
$ws.Range("E3").Value = "  -2.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'499.17"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "'127.93"
$ws.Range("E6").Value = "  -3.85%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").Value = "2.332.06"
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'4.80"
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("E13").Value = "  -1.24%  "
$ws.Range("D14").Value = "2.746.49"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "55.439.37"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "2.344.90"
$ws.Range("E18").Value = "  -2.25%  "
$ws.Range("D19").Value = "'9.86"
$ws.Range("E19").Value = "  -4.12%  "
$ws.Range("D20").Value = "'307.39"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").Value = "'3.99"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'6.14"
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "'65.11"
$ws.Range("E24").Value = "  -3.78%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("D27").Value = "'0.145"
$ws.Range("E27").Value = "  -4.72%  "
$ws.Range("E28").Value = "  -4.55%  "
$ws.Range("D29").Value = "'172.86"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("D30").Value = "'1.62"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").Value = "0.0₃0697"
$ws.Range("E31").Value = "  -4.17%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("D36").Value = "'17.52"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'3.62"
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("D40").Value = "'36.06"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  -5.16%  "
$ws.Range("E42").Value = "  -1.36%  "
$ws.Range("D43").Value = "'126.20"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("D44").Value = "'4.67"
$ws.Range("E44").Value = "  -4.35%  "
$ws.Range("D45").Value = "'0.554"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").Value = "'0.0890"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").Value = "'235.14"
$ws.Range("E47").Value = "  -6.63%  "
$ws.Range("E48").Value = "  -3.09%  "
$ws.Range("E49").Value = "  -3.66%  "
$ws.Range("D50").Value = "'16.54"
$ws.Range("E50").Value = "  -3.29%  "
$ws.Range("D51").Value = "'0.949"
$ws.Range("E51").Value = "  -0.16%  "
